$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain stored as text (matching source data),
# so Excel does not silently coerce numeric-looking strings into floats.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '25.767.92'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.635.68'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '215.43'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '19.88'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '1.638.80'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = '1.861.43'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').Value = '0.0₃0776'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').Value = '63.14'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '25.806.39'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '4.46'
$ws.Range('E20').Value = '  +3.34%  '
$ws.Range('D21').Value = '194.15'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '9.94'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -5.09%  '
$ws.Range('D28').Value = '6.85'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '15.57'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '0.0492'
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('D33').Value = '3.26'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('E34').Value = '  +2.00%  '
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').Value = '0.899'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').Value = '1.113.32'
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').Value = '0.0157'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').Value = '5.58'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '99.19'
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.801'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '0.0₆0110'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('D46').Value = '55.42'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = '2.49'
$ws.Range('E47').Value = '  +11.85%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '7.71'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.417'
$ws.Range('E49').Value = '  -2.79%  '
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.23%  '
